# Generate Report for Handback
# Removes the "cbe7af25-687c-4191-a915-7ed22c306881" handback record (row 3)
# from every worksheet and refreshes the handback timestamps for the
# remaining "b1e3693c-9008-434f-9246-cd16926eb1aa" record.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Hyperlink deletion in this engine only works at the collection level (it
# wipes every hyperlink on the sheet), so remove them all, shrink the grid,
# then recreate the hyperlink(s) that should remain.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/cc9ee3933692a70f6a0b13275a27dbdc0851716e/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
)

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

# Refresh the handoff/handback datetimes for the surviving record before the
# row shift so the writes land on row 2.
$wsZhCn.Range("E2").Value = "2016-03-23 08:54:36"
$wsZhCn.Range("H2").Value = "2016-03-23 08:54:58"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/cc9ee3933692a70f6a0b13275a27dbdc0851716e/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5bdb0691393a4c0441dbf7195855b4d41f4ea70/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5624adc9e3139f84075fcdf6818ef4225f64f894/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e56f238631cef3ac778a5448369854532c20911a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range("E2").Value = "2016-03-23 08:54:40"
$wsDeDe.Range("H2").Value = "2016-03-23 08:55:07"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/cc9ee3933692a70f6a0b13275a27dbdc0851716e/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8328874d47cfbf70fe4fb74f8221fff2c07868e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4e0abbaeaacb4984e498c10ca19792c609b48bd4/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/24ae427b664128de15acf114b90c5414395ba7e9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf",
    "",
    "",
    "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf"
)
